# rnaSample_1660.xlsx -- "fixed harvester column in rnasamples -- holly added
# S.GISH to harvester in bioSamples"
#
# The "harvester" column is column B. Every data row (2-22) had the wrong
# value copied into it (the retrofit label that belongs in rnaPreparer); fix
# it by setting the harvester for each sample to "S.GISH".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B22").Value = "S.GISH"

# Column B visibly narrows/widens to fit the new text once it's typed in.
$ws.Columns.Item(2).ColumnWidth = 8

# Leave the selection on the harvester column, matching the saved view.
$ws.Range("B:B").Select()
